# Apply trade #25 close update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.2
$wsSummary.Range("B4").Value = -0.8100000000000001
$wsSummary.Range("B5").Value = -0.65
$wsSummary.Range("B6").Value = 25
$wsSummary.Range("B7").Value = 8
$wsSummary.Range("B9").Value = 32

# ---- Strategy Status sheet (MarketMaking row, row 4) ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.2
$wsStatus.Range("D4").Value = 25
$wsStatus.Range("E4").Value = -0.8100000000000001
$wsStatus.Range("F4").Value = -0.8
$wsStatus.Range("G4").Value = 32

# ---- Helper to append the new trade #25 row to a trades log sheet ----
function Add-Trade25Row($ws) {
    $ws.Cells.Item(26, 1).Value = 25

    # Force column B to text so the date-like string "2026-02-17" is not
    # auto-converted into a date serial number by the engine.
    $ws.Range("B26").NumberFormat = "@"
    $ws.Cells.Item(26, 2).Value = "2026-02-17"

    $ws.Cells.Item(26, 3).Value = "13:18:45"
    $ws.Cells.Item(26, 4).Value = "MarketMaking"
    $ws.Cells.Item(26, 5).Value = "UP"
    $ws.Cells.Item(26, 6).Value = 0.14
    $ws.Cells.Item(26, 7).Value = 0.17
    $ws.Cells.Item(26, 8).Value = "CLOSED"
    $ws.Cells.Item(26, 9).Value = 21.4286
    $ws.Cells.Item(26, 10).Value = 0.03
    $ws.Cells.Item(26, 11).Value = 99.2
    $ws.Cells.Item(26, 12).Value = 0
    $ws.Cells.Item(26, 13).Value = 0
    $ws.Cells.Item(26, 14).Value = 0.6
    $ws.Cells.Item(26, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(26, 16).Value = "early_exit"
    $ws.Cells.Item(26, 17).Value = 0.13
}

# ---- All Trades sheet ----
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-Trade25Row $wsAllTrades

# ---- MarketMaking sheet ----
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade25Row $wsMarketMaking
